$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Correct the unit prices in column D (price correction due to Google Drive request bug)
$ws.Range("D28").Value = 297.532
$ws.Range("D29").Value = 311.739
$ws.Range("D30").Value = 502.246
$ws.Range("D31").Value = 713.183
$ws.Range("D32").Value = 996.494
$ws.Range("D33").Value = 1278.921
$ws.Range("D34").Value = 1998.314
